# Update the "取得日時" (acquired datetime) column on the active sheet
# from "2025-12-11 18:26:59" to "2025-12-11 18:36:27" for data rows 2-15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-12-11 18:36:27"

for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
